$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "60.856.20"
$ws.Range("E2").Value2 = "  -0.31%  "
$ws.Range("D3").Value2 = "2.401.28"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "561.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  -1.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "141.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "  +1.14%  "
$ws.Range("E7").Value2 = "  -0.34%  "
$ws.Range("E8").Value2 = "  +1.57%  "
$ws.Range("D9").Value2 = "2.408.57"
$ws.Range("E9").Value2 = "  -0.16%  "
$ws.Range("E10").Value2 = "  +0.10%  "
$ws.Range("E11").Value2 = "  -0.43%  "
$ws.Range("E12").Value2 = "  +1.37%  "
$ws.Range("E13").Value2 = "  +1.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "26.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = "  -0.23%  "
$ws.Range("E15").Value2 = "  -1.31%  "
$ws.Range("D16").Value2 = "2.785.40"
$ws.Range("E16").Value2 = "  -3.00%  "
$ws.Range("D17").Value2 = "60.680.48"
$ws.Range("E17").Value2 = "  -0.39%  "
$ws.Range("D18").Value2 = "2.409.18"
$ws.Range("E18").Value2 = "  -0.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "8.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = "  +6.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "10.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "  -0.21%  "
$ws.Range("E21").Value2 = "  +0.02%  "
$ws.Range("E22").Value2 = "  +0.64%  "
$ws.Range("E23").Value2 = "  -0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "0.999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = "  -0.33%  "
$ws.Range("E25").Value2 = "  -2.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "64.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = "  -0.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "571.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value2 = "  -2.74%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "8.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = "  -4.74%  "
$ws.Range("D30").Value2 = "0.0₃0940"
$ws.Range("E30").Value2 = "  -0.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "8.08"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value2 = "  +2.04%  "
$ws.Range("E32").Value2 = "  -2.09%  "
$ws.Range("E33").Value2 = "  -2.50%  "
$ws.Range("E34").Value2 = "  +0.00%  "
$ws.Range("E35").Value2 = "  -0.54%  "
$ws.Range("E36").Value2 = "  +3.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "153.43"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value2 = "  +0.95%  "
$ws.Range("E38").Value2 = "  +0.34%  "
$ws.Range("E39").Value2 = "  -1.29%  "
$ws.Range("E40").Value2 = "  -0.09%  "
$ws.Range("E41").Value2 = "  +0.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "2.56"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value2 = "  +7.94%  "
$ws.Range("E43").Value2 = "  -0.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "41.95"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = "  +1.70%  "
$ws.Range("E45").Value2 = "  -0.71%  "
$ws.Range("D46").Value2 = "0.0₆0278"
$ws.Range("E46").Value2 = "  -5.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "142.14"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = "  -0.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "3.51"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "  -0.32%  "
$ws.Range("E49").Value2 = "  -0.68%  "
$ws.Range("E50").Value2 = "  +0.27%  "
$ws.Range("E51").Value2 = "  -1.88%  "
